$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.887.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.857.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.858.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.43%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.483.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.846.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.969.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -4.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.008.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("E33").Value = "  -3.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.832.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.46%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "428.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "47.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
